$wb = $excel.ActiveWorkbook

# ============================================================
# Step 1: Insert new worksheet "2022-Q3" immediately after "总计"
#         (i.e. before the current sheet at position 2, "2022-Q2")
# ============================================================
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$newSheet.Name = "2022-Q3"

$ws = $wb.Worksheets.Item(2)

# Match sheetPr/outlinePr settings used by the other quarter sheets
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# ---- Header row ----
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# ---- Data rows (0-based index in col A, fund code/size/position/ratio/value kept as text) ----
    $ws.Cells.Item(2, 1).Value = 0
    $ws.Cells.Item(2, 2).Value = "'159941"
    $ws.Cells.Item(2, 3).Value = "广发纳斯达克100ETF（QDII）"
    $ws.Cells.Item(2, 4).Value = "'106.15"
    $ws.Cells.Item(2, 5).Value = "'91.14"
    $ws.Cells.Item(2, 6).Value = "'3.12"
    $ws.Cells.Item(2, 7).Value = "'3.3119"
    $ws.Cells.Item(2, 8).Value = 6
    $ws.Cells.Item(3, 1).Value = 1
    $ws.Cells.Item(3, 2).Value = "'513100"
    $ws.Cells.Item(3, 3).Value = "国泰纳斯达克100（QDII-ETF）"
    $ws.Cells.Item(3, 4).Value = "'46.54"
    $ws.Cells.Item(3, 5).Value = "'91.35"
    $ws.Cells.Item(3, 6).Value = "'3.12"
    $ws.Cells.Item(3, 7).Value = "'1.4520"
    $ws.Cells.Item(3, 8).Value = 6
    $ws.Cells.Item(4, 1).Value = 2
    $ws.Cells.Item(4, 2).Value = "'513500"
    $ws.Cells.Item(4, 3).Value = "博时标普500ETF（QDII）"
    $ws.Cells.Item(4, 4).Value = "'71.37"
    $ws.Cells.Item(4, 5).Value = "'96.44"
    $ws.Cells.Item(4, 6).Value = "'1.84"
    $ws.Cells.Item(4, 7).Value = "'1.3132"
    $ws.Cells.Item(4, 8).Value = 5
    $ws.Cells.Item(5, 1).Value = 3
    $ws.Cells.Item(5, 2).Value = "'040047"
    $ws.Cells.Item(5, 3).Value = "华安纳斯达克100指数（QDII）美元现钞A"
    $ws.Cells.Item(5, 4).Value = "'24.52"
    $ws.Cells.Item(5, 5).Value = "'92.09"
    $ws.Cells.Item(5, 6).Value = "'3.16"
    $ws.Cells.Item(5, 7).Value = "'0.7748"
    $ws.Cells.Item(5, 8).Value = 6
    $ws.Cells.Item(6, 1).Value = 4
    $ws.Cells.Item(6, 2).Value = "'040048"
    $ws.Cells.Item(6, 3).Value = "华安纳斯达克100指数（QDII）美元现汇A"
    $ws.Cells.Item(6, 4).Value = "'24.52"
    $ws.Cells.Item(6, 5).Value = "'92.09"
    $ws.Cells.Item(6, 6).Value = "'3.16"
    $ws.Cells.Item(6, 7).Value = "'0.7748"
    $ws.Cells.Item(6, 8).Value = 6
    $ws.Cells.Item(7, 1).Value = 5
    $ws.Cells.Item(7, 2).Value = "'040046"
    $ws.Cells.Item(7, 3).Value = "华安纳斯达克100指数（QDII）人民币A"
    $ws.Cells.Item(7, 4).Value = "'22.21"
    $ws.Cells.Item(7, 5).Value = "'92.09"
    $ws.Cells.Item(7, 6).Value = "'3.16"
    $ws.Cells.Item(7, 7).Value = "'0.7018"
    $ws.Cells.Item(7, 8).Value = 6
    $ws.Cells.Item(8, 1).Value = 6
    $ws.Cells.Item(8, 2).Value = "'000043"
    $ws.Cells.Item(8, 3).Value = "嘉实美国成长股票（QDII）人民币"
    $ws.Cells.Item(8, 4).Value = "'12.41"
    $ws.Cells.Item(8, 5).Value = "'92.80"
    $ws.Cells.Item(8, 6).Value = "'5.46"
    $ws.Cells.Item(8, 7).Value = "'0.6776"
    $ws.Cells.Item(8, 8).Value = 3
    $ws.Cells.Item(9, 1).Value = 7
    $ws.Cells.Item(9, 2).Value = "'000044"
    $ws.Cells.Item(9, 3).Value = "嘉实美国成长股票（QDII）美元现汇"
    $ws.Cells.Item(9, 4).Value = "'12.41"
    $ws.Cells.Item(9, 5).Value = "'92.80"
    $ws.Cells.Item(9, 6).Value = "'5.46"
    $ws.Cells.Item(9, 7).Value = "'0.6776"
    $ws.Cells.Item(9, 8).Value = 3
    $ws.Cells.Item(10, 1).Value = 8
    $ws.Cells.Item(10, 2).Value = "'001668"
    $ws.Cells.Item(10, 3).Value = "汇添富全球移动互联灵活配置混合（QDII）A"
    $ws.Cells.Item(10, 4).Value = "'12.06"
    $ws.Cells.Item(10, 5).Value = "'90.88"
    $ws.Cells.Item(10, 6).Value = "'3.86"
    $ws.Cells.Item(10, 7).Value = "'0.4655"
    $ws.Cells.Item(10, 8).Value = 3
    $ws.Cells.Item(11, 1).Value = 9
    $ws.Cells.Item(11, 2).Value = "'160213"
    $ws.Cells.Item(11, 3).Value = "国泰纳斯达克100指数（QDII）"
    $ws.Cells.Item(11, 4).Value = "'15.14"
    $ws.Cells.Item(11, 5).Value = "'85.81"
    $ws.Cells.Item(11, 6).Value = "'3.00"
    $ws.Cells.Item(11, 7).Value = "'0.4542"
    $ws.Cells.Item(11, 8).Value = 6
    $ws.Cells.Item(12, 1).Value = 10
    $ws.Cells.Item(12, 2).Value = "'000834"
    $ws.Cells.Item(12, 3).Value = "大成纳斯达克100指数（QDII）"
    $ws.Cells.Item(12, 4).Value = "'14.15"
    $ws.Cells.Item(12, 5).Value = "'85.22"
    $ws.Cells.Item(12, 6).Value = "'2.92"
    $ws.Cells.Item(12, 7).Value = "'0.4132"
    $ws.Cells.Item(12, 8).Value = 6
    $ws.Cells.Item(13, 1).Value = 11
    $ws.Cells.Item(13, 2).Value = "'513300"
    $ws.Cells.Item(13, 3).Value = "华夏纳斯达克100ETF（QDII）"
    $ws.Cells.Item(13, 4).Value = "'11.08"
    $ws.Cells.Item(13, 5).Value = "'97.32"
    $ws.Cells.Item(13, 6).Value = "'3.34"
    $ws.Cells.Item(13, 7).Value = "'0.3701"
    $ws.Cells.Item(13, 8).Value = 1
    $ws.Cells.Item(14, 1).Value = 12
    $ws.Cells.Item(14, 2).Value = "'003722"
    $ws.Cells.Item(14, 3).Value = "易方达纳斯达克100指数美元（QDII-LOF）A"
    $ws.Cells.Item(14, 4).Value = "'7.72"
    $ws.Cells.Item(14, 5).Value = "'90.67"
    $ws.Cells.Item(14, 6).Value = "'3.10"
    $ws.Cells.Item(14, 7).Value = "'0.2393"
    $ws.Cells.Item(14, 8).Value = 6
    $ws.Cells.Item(15, 1).Value = 13
    $ws.Cells.Item(15, 2).Value = "'161130"
    $ws.Cells.Item(15, 3).Value = "易方达纳斯达克100指数人民币（QDII-LOF）"
    $ws.Cells.Item(15, 4).Value = "'7.72"
    $ws.Cells.Item(15, 5).Value = "'90.67"
    $ws.Cells.Item(15, 6).Value = "'3.10"
    $ws.Cells.Item(15, 7).Value = "'0.2393"
    $ws.Cells.Item(15, 8).Value = 6
    $ws.Cells.Item(16, 1).Value = 14
    $ws.Cells.Item(16, 2).Value = "'100055"
    $ws.Cells.Item(16, 3).Value = "富国全球科技互联网股票（QDII）"
    $ws.Cells.Item(16, 4).Value = "'3.95"
    $ws.Cells.Item(16, 5).Value = "'86.97"
    $ws.Cells.Item(16, 6).Value = "'4.64"
    $ws.Cells.Item(16, 7).Value = "'0.1833"
    $ws.Cells.Item(16, 8).Value = 6
    $ws.Cells.Item(17, 1).Value = 15
    $ws.Cells.Item(17, 2).Value = "'161125"
    $ws.Cells.Item(17, 3).Value = "易方达标普500指数（QDII-LOF）人民币"
    $ws.Cells.Item(17, 4).Value = "'4.74"
    $ws.Cells.Item(17, 5).Value = "'90.72"
    $ws.Cells.Item(17, 6).Value = "'1.74"
    $ws.Cells.Item(17, 7).Value = "'0.0825"
    $ws.Cells.Item(17, 8).Value = 5
    $ws.Cells.Item(18, 1).Value = 16
    $ws.Cells.Item(18, 2).Value = "'012860"
    $ws.Cells.Item(18, 3).Value = "易方达标普500指数（QDII-LOF）人民币 C"
    $ws.Cells.Item(18, 4).Value = "'4.74"
    $ws.Cells.Item(18, 5).Value = "'90.72"
    $ws.Cells.Item(18, 6).Value = "'1.74"
    $ws.Cells.Item(18, 7).Value = "'0.0825"
    $ws.Cells.Item(18, 8).Value = 5
    $ws.Cells.Item(19, 1).Value = 17
    $ws.Cells.Item(19, 2).Value = "'160644"
    $ws.Cells.Item(19, 3).Value = "鹏华香港美国互联网股票（LOF）人民币"
    $ws.Cells.Item(19, 4).Value = "'1.23"
    $ws.Cells.Item(19, 5).Value = "'83.13"
    $ws.Cells.Item(19, 6).Value = "'6.63"
    $ws.Cells.Item(19, 7).Value = "'0.0815"
    $ws.Cells.Item(19, 8).Value = 3
    $ws.Cells.Item(20, 1).Value = 18
    $ws.Cells.Item(20, 2).Value = "'006792"
    $ws.Cells.Item(20, 3).Value = "鹏华香港美国互联网股票（LOF）美元现汇"
    $ws.Cells.Item(20, 4).Value = "'1.23"
    $ws.Cells.Item(20, 5).Value = "'83.13"
    $ws.Cells.Item(20, 6).Value = "'6.63"
    $ws.Cells.Item(20, 7).Value = "'0.0815"
    $ws.Cells.Item(20, 8).Value = 3
    $ws.Cells.Item(21, 1).Value = 19
    $ws.Cells.Item(21, 2).Value = "'003718"
    $ws.Cells.Item(21, 3).Value = "易方达标普500指数（QDII-LOF）美元A"
    $ws.Cells.Item(21, 4).Value = "'4.66"
    $ws.Cells.Item(21, 5).Value = "'90.72"
    $ws.Cells.Item(21, 6).Value = "'1.74"
    $ws.Cells.Item(21, 7).Value = "'0.0811"
    $ws.Cells.Item(21, 8).Value = 5
    $ws.Cells.Item(22, 1).Value = 20
    $ws.Cells.Item(22, 2).Value = "'014978"
    $ws.Cells.Item(22, 3).Value = "华安纳斯达克100指数（QDII）人民币C"
    $ws.Cells.Item(22, 4).Value = "'2.31"
    $ws.Cells.Item(22, 5).Value = "'92.09"
    $ws.Cells.Item(22, 6).Value = "'3.16"
    $ws.Cells.Item(22, 7).Value = "'0.0730"
    $ws.Cells.Item(22, 8).Value = 6
    $ws.Cells.Item(23, 1).Value = 21
    $ws.Cells.Item(23, 2).Value = "'159632"
    $ws.Cells.Item(23, 3).Value = "华安纳斯达克100ETF（QDII）"
    $ws.Cells.Item(23, 4).Value = "'1.51"
    $ws.Cells.Item(23, 5).Value = "'89.05"
    $ws.Cells.Item(23, 6).Value = "'3.07"
    $ws.Cells.Item(23, 7).Value = "'0.0464"
    $ws.Cells.Item(23, 8).Value = 6
    $ws.Cells.Item(24, 1).Value = 22
    $ws.Cells.Item(24, 2).Value = "'013329"
    $ws.Cells.Item(24, 3).Value = "嘉实全球价值股票（QDII）美元现汇"
    $ws.Cells.Item(24, 4).Value = "'1.68"
    $ws.Cells.Item(24, 5).Value = "'90.63"
    $ws.Cells.Item(24, 6).Value = "'1.95"
    $ws.Cells.Item(24, 7).Value = "'0.0328"
    $ws.Cells.Item(24, 8).Value = 4
    $ws.Cells.Item(25, 1).Value = 23
    $ws.Cells.Item(25, 2).Value = "'013328"
    $ws.Cells.Item(25, 3).Value = "嘉实全球价值股票（QDII）人民币"
    $ws.Cells.Item(25, 4).Value = "'1.68"
    $ws.Cells.Item(25, 5).Value = "'90.63"
    $ws.Cells.Item(25, 6).Value = "'1.95"
    $ws.Cells.Item(25, 7).Value = "'0.0328"
    $ws.Cells.Item(25, 8).Value = 4
    $ws.Cells.Item(26, 1).Value = 24
    $ws.Cells.Item(26, 2).Value = "'159612"
    $ws.Cells.Item(26, 3).Value = "国泰标普500ETF（QDII）"
    $ws.Cells.Item(26, 4).Value = "'0.55"
    $ws.Cells.Item(26, 5).Value = "'91.40"
    $ws.Cells.Item(26, 6).Value = "'1.73"
    $ws.Cells.Item(26, 7).Value = "'0.0095"
    $ws.Cells.Item(26, 8).Value = 5
    $ws.Cells.Item(27, 1).Value = 25
    $ws.Cells.Item(27, 2).Value = "'012871"
    $ws.Cells.Item(27, 3).Value = "易方达纳斯达克100指数美元（QDII-LOF）C"
    $ws.Cells.Item(27, 4).Value = "'0.18"
    $ws.Cells.Item(27, 5).Value = "'90.67"
    $ws.Cells.Item(27, 6).Value = "'3.10"
    $ws.Cells.Item(27, 7).Value = "'0.0056"
    $ws.Cells.Item(27, 8).Value = 6
    $ws.Cells.Item(28, 1).Value = 26
    $ws.Cells.Item(28, 2).Value = "'012870"
    $ws.Cells.Item(28, 3).Value = "易方达纳斯达克100指数人民币（QDII-LOF）C"
    $ws.Cells.Item(28, 4).Value = "'0.18"
    $ws.Cells.Item(28, 5).Value = "'90.67"
    $ws.Cells.Item(28, 6).Value = "'3.10"
    $ws.Cells.Item(28, 7).Value = "'0.0056"
    $ws.Cells.Item(28, 8).Value = 6
    $ws.Cells.Item(29, 1).Value = 27
    $ws.Cells.Item(29, 2).Value = "'015203"
    $ws.Cells.Item(29, 3).Value = "汇添富全球移动互联灵活配置混合（QDII）D"
    $ws.Cells.Item(29, 4).Value = "'0.04"
    $ws.Cells.Item(29, 5).Value = "'90.88"
    $ws.Cells.Item(29, 6).Value = "'3.86"
    $ws.Cells.Item(29, 7).Value = "'0.0015"
    $ws.Cells.Item(29, 8).Value = 3
    $ws.Cells.Item(30, 1).Value = 28
    $ws.Cells.Item(30, 2).Value = "'012861"
    $ws.Cells.Item(30, 3).Value = "易方达标普500指数（QDII-LOF）美元 C"
    $ws.Cells.Item(30, 4).Value = "'0.08"
    $ws.Cells.Item(30, 5).Value = "'90.72"
    $ws.Cells.Item(30, 6).Value = "'1.74"
    $ws.Cells.Item(30, 7).Value = "'0.0014"
    $ws.Cells.Item(30, 8).Value = 5
    $ws.Cells.Item(31, 1).Value = 29
    $ws.Cells.Item(31, 2).Value = "'015202"
    $ws.Cells.Item(31, 3).Value = "汇添富全球移动互联灵活配置混合（QDII）C"
    $ws.Cells.Item(31, 4).Value = "'0.01"
    $ws.Cells.Item(31, 5).Value = "'90.88"
    $ws.Cells.Item(31, 6).Value = "'3.86"
    $ws.Cells.Item(31, 7).Value = "'0.0004"
    $ws.Cells.Item(31, 8).Value = 3

# Clear the automatic "text" number-format style that got applied when
# entering apostrophe-prefixed numeric-looking strings, so these cells end
# up with no explicit style (matching a freshly written text cell).
$ws.Range("B2:G31").Style = "Normal"

# Apply the correct header style (s=2) and index-column style (s=2) by
# copying formatting from the equivalent cells of a sibling quarter sheet
# (now shifted to position 3, e.g. "2022-Q2") which already uses them.
$srcWs = $wb.Worksheets.Item(3)
$srcWs.Range("B1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$srcWs.Range("A2").Copy()
$ws.Range("A2:A31").PasteSpecial(-4122)

# ============================================================
# Step 2: Insert a new row into "总计" (summary) sheet for 2022-Q3
# ============================================================
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()

# Clear formatting picked up from the row above during insert
$summary.Range("B2:D2").Style = "Normal"

# Apply index-column style (s=2) to the new A2 cell
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

# New summary row values
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 30
$summary.Cells.Item(2, 4).Value = 12.67

# Re-sequence the 0-based index column for the rows that shifted down
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5

# Ensure the summary sheet ("总计") remains the active tab, matching original view state
$wb.Worksheets.Item(1).Activate()
